$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pirmas lapas")

$ws.Range("A2").Value = "Cia dar kazkas prasyta Windows"
$ws.Range("A5").Select()
